# All command done! Test one final time and check error handling with
# appropriate message. Also add comments for dtar

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old "display" scratch block that lived at rows 36-38 ---
# (its surviving values get rewritten further down at their new home)
$ws.Range("A36").ClearContents()
$ws.Range("B36").ClearContents()
$ws.Range("B37").ClearContents()
$ws.Range("B38").ClearContents()

# --- Update condition text for the dfile block (server -> servers) ---
$ws.Range("D27").Value = "in servers"

# --- "display" test matrix moved down to rows 45-50 with new columns ---
$ws.Range("A45").Value = "display ~/smain/test"
$ws.Range("A46").Value = "display ~/smain/test2"
$ws.Range("C45").Value = "test1 in main"
$ws.Range("C46").Value = "test2 in pdf and text"
$ws.Range("C47").Value = "test 3 in main and pdf"
$ws.Range("C48").Value = "test 4 in main and text"
$ws.Range("A47").Value = "display ~/smain/test3"
$ws.Range("A48").Value = "display ~/smain/test4"
$ws.Range("B45").Value = "display ~/smain/test5"
$ws.Range("C49").Value = "test5 in none"
$ws.Range("B48").Value = "display ~/smain/test6"
$ws.Range("C50").Value = "test 6 in all but empty"
$ws.Range("B46").Value = "display ~/spdf/"
$ws.Range("B47").Value = "display ~/stxt/"

# --- New "dtar" comments block (rows 35-37) ---
$ws.Range("A35").Value = "dtar .c"
$ws.Range("A36").Value = "dtar .pdf"
$ws.Range("A37").Value = "dtar .txt"
$ws.Range("B35").Value = "dtar .pptx"
$ws.Range("B36").Value = "dtar .xlxs"

# --- New "ufile" quick reference rows under the first block ---
$ws.Range("B9").Value = "ufile ~/smain/"
$ws.Range("B10").Value = "ufile"

# --- View state: zoom in and move the selection/scroll position ---
$excel.ActiveWindow.Zoom = 115
$ws.Range("B10").Select()
